$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing header M1
$ws.Range("M1").Value = "Member name"

# Add new headers N1:Q1 (same header style as M1, i.e. style index 1)
$ws.Range("N1").Value = "Article name"
$ws.Range("O1").Value = "Page number"
$ws.Range("P1").Value = "Tag"
$ws.Range("Q1").Value = "Category"

$ws.Range("M1").Copy()
$ws.Range("N1:Q1").PasteSpecial(-4122)

# Column widths for the new columns N and O (closest achievable values to the
# target 12.14 / 12.91 "characters" widths under this engine's column-width
# quantization, which snaps to 1/6-character increments)
$ws.Columns.Item(14).ColumnWidth = 11.3
$ws.Columns.Item(15).ColumnWidth = 12.0

# Update the view: scroll to show column I at top-left, select Q2
$ws.Application.ActiveWindow.ScrollColumn = 9
$ws.Range("Q2").Select() | Out-Null
